$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.817.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.115.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.65%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5319"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.00%  "

$ws.Range("E8").Value = "  +6.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09009"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.74%  "

$ws.Range("E10").Value = "  +8.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.180"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.108.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.762"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.844"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.97%  "

$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001130"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.68%  "

$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("E22").Value = "  +3.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.870.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.357.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.264"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.581"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.175"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.84%  "

$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.227"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.017"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.557"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +18.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02620"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.80%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.03%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.539"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06750"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.503"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2284"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6850"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.249"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.08%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.17%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6451"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.52%  "

$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.231"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.661"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("E49").Value = "  +4.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.45%  "
